# feat: add 2022-Q1 data
#
# - insert a new worksheet "2022-Q1" between "2021-Q3" and "总计",
#   populated with the new fund-holding rows
# - insert a new top row in "总计" summarising the new quarter and
#   push the existing "2021-Q3" summary row down

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as genuine TEXT (so numeric-looking
# strings like "011815" or "4.64" are not silently coerced into
# numbers), while leaving the cell's style untouched (falls back to the
# sheet default "no style" look, same as every other plain data cell).
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet in place of the old "总计" sheet,
#    and re-create a fresh "总计" sheet after it (this keeps sheetIds
#    sequential: 2021-Q3=1, 2022-Q1=2, 总计=3 -- matching how a sheet
#    that already existed keeps a lower id than one newly added).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Duplicate "总计" so the copy can become the new, empty "总计" sheet at
# the end of the book (this preserves the sheet's pane/format setup).
$totalSheet.Copy($null, $totalSheet)

# The original sheet object now becomes "2022-Q1".
$qSheet = $totalSheet
$qSheet.Name = "2022-Q1"

# The freshly made copy (now the last sheet) becomes the new "总计".
$newTotalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newTotalSheet.Name = "总计"

# ---------------------------------------------------------------------
# 2) Populate "2022-Q1" with the fund-holding table.
# ---------------------------------------------------------------------

# Headers (row 1) -- reuse the existing bold/bordered header style that
# was already on B1:D1 by copying its format onto the new header cells.
$qSheet.Range("B1").Copy() | Out-Null
$qSheet.Range("E1:H1").PasteSpecial(-4122)
$qSheet.Application.CutCopyMode = $false

$qSheet.Range("B1").Value = "基金代码"
$qSheet.Range("C1").Value = "基金名称"
$qSheet.Range("D1").Value = "基金规模"
$qSheet.Range("E1").Value = "股票总仓位"
$qSheet.Range("F1").Value = "仓位占比"
$qSheet.Range("G1").Value = "持有市值(亿元)"
$qSheet.Range("H1").Value = "仓位排名"

# Row 2 -- first fund. A2 keeps the existing "index" style (s=2) and is
# a plain number, exactly like the row this sheet's data replaces.
$qSheet.Range("A2").Value = 0
Set-TextValue $qSheet.Range("B2") "011815"
$qSheet.Range("C2").Value = "恒越优势精选混合型发起式证券投资基金"
Set-TextValue $qSheet.Range("D2") "4.64"
Set-TextValue $qSheet.Range("E2") "92.44"
Set-TextValue $qSheet.Range("F2") "3.95"
Set-TextValue $qSheet.Range("G2") "0.1833"
$qSheet.Range("H2").Value = 6

# Row 3 -- second fund. A3 needs the same index style as A2.
$qSheet.Range("A2").Copy() | Out-Null
$qSheet.Range("A3").PasteSpecial(-4122)
$qSheet.Application.CutCopyMode = $false

$qSheet.Range("A3").Value = 1
Set-TextValue $qSheet.Range("B3") "013028"
$qSheet.Range("C3").Value = "恒越品质生活混合"
Set-TextValue $qSheet.Range("D3") "2.03"
Set-TextValue $qSheet.Range("E3") "92.89"
Set-TextValue $qSheet.Range("F3") "3.81"
Set-TextValue $qSheet.Range("G3") "0.0773"
$qSheet.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 3) Update "总计": push the existing 2021-Q3 row down one row and add
#    the new 2022-Q1 summary row above it.
# ---------------------------------------------------------------------
$newTotalSheet.Rows(2).Insert()
$newTotalSheet.Range("A2:D2").ClearFormats()

# A2 should carry the same "index" style that A3 (the shifted row) has,
# and the shifted row's own index needs to move from 0 to 1.
$newTotalSheet.Range("A3").Copy() | Out-Null
$newTotalSheet.Range("A2").PasteSpecial(-4122)
$newTotalSheet.Application.CutCopyMode = $false

$newTotalSheet.Range("A2").Value = 0
$newTotalSheet.Range("B2").Value = "2022-Q1"
$newTotalSheet.Range("C2").Value = 2
$newTotalSheet.Range("D2").Value = 0.26

$newTotalSheet.Range("A3").Value = 1
